# Update matching, add validations
#
# The European Parliament (EP) election rows are being extended with two
# earlier elections (EP 2004, EP 2009) in addition to a fuller EP 2014
# entry, pushing all the following rows (Municipal ..., City districts ...,
# Senate By-elections) down by three rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows right above the old row 18 ("Municipal 1994") so that
# everything from there on (Municipal/City districts/Senate) shifts down by
# three rows without needing to be rewritten.
$ws.Rows("18:20").Insert()

# Row 16: EP 2019 -> EP 2004 (new data)
$ws.Range("A16").Value = "EP 2004"
$ws.Range("B16").Value = 806
$ws.Range("C16").Value = 203
$ws.Range("D16").Value = 24
$ws.Range("E16").Value = 5
$ws.Range("F16").Value = 3

# Row 17: EP 2024 -> EP 2009 (new data)
$ws.Range("A17").Value = "EP 2009"
$ws.Range("B17").Value = 708
$ws.Range("C17").Value = 200
$ws.Range("D17").Value = 22
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = 3.1

# Row 18 (newly inserted): EP 2014
$ws.Range("A18").Value = "EP 2014"
$ws.Range("B18").Value = 849
$ws.Range("C18").Value = 227
$ws.Range("D18").Value = 21
$ws.Range("E18").Value = 5
$ws.Range("F18").Value = 2.5

# Row 19 (newly inserted): EP 2019 (restored, same values as the old row 16)
$ws.Range("A19").Value = "EP 2019"
$ws.Range("B19").Value = 841
$ws.Range("C19").Value = 201
$ws.Range("D19").Value = 21
$ws.Range("E19").Value = 7
$ws.Range("F19").Value = 2.5

# Row 20 (newly inserted): EP 2024 (restored, same values as the old row 17)
$ws.Range("A20").Value = "EP 2024"
$ws.Range("B20").Value = 674
$ws.Range("C20").Value = 245
$ws.Range("D20").Value = 21
$ws.Range("E20").Value = 8
$ws.Range("F20").Value = 3.1

# Rows 21-37 (Municipal 1994 .. Senate By-elections) keep their original
# values - they were simply shifted down by the insert above.
